# Update the "Final Effluent Quality" sheet: rows 2-7 (columns A-D) need to be
# rewritten so that the (label, value, target, max) tuples end up in a new
# row order, as produced by the upstream refactor.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Final Effluent Quality")

# Final (row -> values) state taken from the target workbook.
# Column A values are the shared-string labels (looked up via sharedStrings
# index in the diff); we just set the text directly here.
$rows = @(
    @{ Row = 2; Label = "Target_Effluent_TSS (mg/L)";   Value = 9.792375320235159;  Target = 70; Max = 200 },
    @{ Row = 3; Label = "Target_Effluent_S_NO3 (mg/L)"; Value = 7.437532810673214;  Target = 14; Max = 30 },
    @{ Row = 4; Label = "Target_Effluent_BOD (mg/L)";   Value = 28.64803647229185; Target = 30; Max = 100 },
    @{ Row = 5; Label = "Target_Effluent_COD (mg/L)";   Value = 29.26060168082566; Target = 60; Max = 300 },
    @{ Row = 6; Label = "Target_Effluent_S_NH4 (mg/L)"; Value = 8.049427243199435; Target = 2;  Max = 9 },
    @{ Row = 7; Label = "Target_Effluent_TP (mg/L)";    Value = 8.777835146452022; Target = 1;  Max = 10 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Label
    $ws.Cells.Item($r.Row, 2).Value = $r.Value
    $ws.Cells.Item($r.Row, 3).Value = $r.Target
    $ws.Cells.Item($r.Row, 4).Value = $r.Max
}
